$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 801 (shifts existing rows 801-842 down to 802-843)
$ws.Rows("801:801").Insert()

# Populate the newly inserted row 801 with the new entry.
# The date-like text must stay literal text (not be auto-converted to a
# date serial by Excel's type inference), so it is entered the same way
# a user would force text in the Excel UI: a leading apostrophe.
$ws.Range("A801").Value = "'2026/02/16"
$ws.Range("B801").Value = "月"
$ws.Range("C801").Value = 0
$ws.Range("D801").Value = 35
